# Scheduled-runner update: refresh cached market-price figures (columns
# H:N) on several Leve profit rows across the per-job sheets. Values are
# literal numbers (no formulas in this workbook), so each changed cell is
# written directly via Range.Value; the couple of cells that gain or lose
# a value entirely are handled with an explicit Value assignment / ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 43481796
$ws.Range("I76").Value = 71432410
$ws.Range("J76").Value = 3066.6667
$ws.Range("K76").Value = 71432410
$ws.Range("L76").Value = 3066.6667
$ws.Range("M76").Value = -71432095
$ws.Range("N76").Value = -3696.6667

$ws.Range("H79").Value = 43481796
$ws.Range("I79").Value = 71432410
$ws.Range("J79").Value = 3066.6667
$ws.Range("K79").Value = 71432410
$ws.Range("L79").Value = 3066.6667
$ws.Range("M79").Value = -71431318
$ws.Range("N79").Value = -5250.6667

$ws.Range("H132").Value = 2962120
$ws.Range("I132").Value = 598581.0600000001
$ws.Range("J132").Value = 27779278
$ws.Range("K132").Value = 1795743.18
$ws.Range("L132").Value = 83337834
$ws.Range("M132").Value = -1793213.18
$ws.Range("N132").Value = -83342894

$ws.Range("H137").Value = 14716266
$ws.Range("J137").Value = 27196156
$ws.Range("L137").Value = 81588468
$ws.Range("N137").Value = -81593568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1906.6
$ws.Range("I63").Value = 1920.75
$ws.Range("J63").Value = 1850
$ws.Range("K63").Value = 1920.75
$ws.Range("L63").Value = 1850
$ws.Range("M63").Value = -1234.75
$ws.Range("N63").Value = -3222

$ws.Range("H66").Value = 1906.6
$ws.Range("I66").Value = 1920.75
$ws.Range("J66").Value = 1850
$ws.Range("K66").Value = 9603.75
$ws.Range("L66").Value = 9250
$ws.Range("M66").Value = -6171.75
$ws.Range("N66").Value = -16114

$ws.Range("H74").Value = 29488568
$ws.Range("I74").Value = 27027694
$ws.Range("J74").Value = 35558720
$ws.Range("K74").Value = 27027694
$ws.Range("L74").Value = 35558720
$ws.Range("M74").Value = -27026820
$ws.Range("N74").Value = -35560468

$ws.Range("H77").Value = 29488568
$ws.Range("I77").Value = 27027694
$ws.Range("J77").Value = 35558720
$ws.Range("K77").Value = 135138470
$ws.Range("L77").Value = 177793600
$ws.Range("M77").Value = -135134102
$ws.Range("N77").Value = -177802336

$ws.Range("H88").Value = 3432.5715
$ws.Range("I88").Value = 1750
$ws.Range("J88").Value = 4105.6
$ws.Range("K88").Value = 1750
$ws.Range("L88").Value = 4105.6
$ws.Range("M88").Value = -1344
$ws.Range("N88").Value = -4917.6

$ws.Range("H91").Value = 3432.5715
$ws.Range("I91").Value = 1750
$ws.Range("J91").Value = 4105.6
$ws.Range("K91").Value = 1750
$ws.Range("L91").Value = 4105.6
$ws.Range("M91").Value = -346
$ws.Range("N91").Value = -6913.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 527.1905
$ws.Range("I22").Value = 527.1905
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 527.1905
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = -354.1905
$ws.Range("M22").ClearContents()

$ws.Range("H105").Value = 1699.2858
$ws.Range("I105").Value = 1680.909
$ws.Range("J105").Value = 1766.6666
$ws.Range("K105").Value = 1680.909
$ws.Range("L105").Value = 1766.6666
$ws.Range("M105").Value = 66.09099999999989
$ws.Range("N105").Value = -5260.6666

$ws.Range("H134").Value = 10821266
$ws.Range("I134").Value = 11793268
$ws.Range("J134").Value = 129250
$ws.Range("K134").Value = 35379804
$ws.Range("L134").Value = 387750
$ws.Range("M134").Value = -35377269
$ws.Range("N134").Value = -392820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1941018.6
$ws.Range("I31").Value = 2689390.2
$ws.Range("J31").Value = 7724.9165
$ws.Range("K31").Value = 2689390.2
$ws.Range("L31").Value = 7724.9165
$ws.Range("M31").Value = -2689095.2
$ws.Range("N31").Value = -8314.916499999999

$ws.Range("H34").Value = 1941018.6
$ws.Range("I34").Value = 2689390.2
$ws.Range("J34").Value = 7724.9165
$ws.Range("K34").Value = 2689390.2
$ws.Range("L34").Value = 7724.9165
$ws.Range("M34").Value = -2689188.2
$ws.Range("N34").Value = -8128.9165

$ws.Range("H58").Value = 930632.3
$ws.Range("I58").Value = 3247.641
$ws.Range("J58").Value = 4547432.5
$ws.Range("K58").Value = 3247.641
$ws.Range("L58").Value = 4547432.5
$ws.Range("M58").Value = -3044.641
$ws.Range("N58").Value = -4547838.5

$ws.Range("H62").Value = 2706.1052
$ws.Range("I62").Value = 2253.3333
$ws.Range("J62").Value = 3482.2856
$ws.Range("K62").Value = 2253.3333
$ws.Range("L62").Value = 3482.2856
$ws.Range("M62").Value = -1629.3333
$ws.Range("N62").Value = -4730.2856

$ws.Range("H65").Value = 2706.1052
$ws.Range("I65").Value = 2253.3333
$ws.Range("J65").Value = 3482.2856
$ws.Range("K65").Value = 11266.6665
$ws.Range("L65").Value = 17411.428
$ws.Range("M65").Value = -8146.666499999999
$ws.Range("N65").Value = -23651.428

$ws.Range("H136").Value = 930632.3
$ws.Range("I136").Value = 3247.641
$ws.Range("J136").Value = 4547432.5
$ws.Range("K136").Value = 9742.923000000001
$ws.Range("L136").Value = 13642297.5
$ws.Range("M136").Value = -7192.923000000001
$ws.Range("N136").Value = -13647397.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1489356.6
$ws.Range("I5").Value = 680.65
$ws.Range("J5").Value = 2842698.5
$ws.Range("K5").Value = 2041.95
$ws.Range("L5").Value = 8528095.5
$ws.Range("M5").Value = -1929.95
$ws.Range("N5").Value = -8528319.5

$ws.Range("H128").Value = 125500
$ws.Range("I128").Value = 125500
$ws.Range("K128").Value = 376500
$ws.Range("M128").Value = -371520

$ws.Range("H135").Value = 1489356.6
$ws.Range("I135").Value = 680.65
$ws.Range("J135").Value = 2842698.5
$ws.Range("K135").Value = 6125.849999999999
$ws.Range("L135").Value = 25584286.5
$ws.Range("M135").Value = -3590.849999999999
$ws.Range("N135").Value = -25589356.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 15406
$ws.Range("J40").Value = 15406
$ws.Range("L40").Value = 15406
$ws.Range("N40").Value = -15708

$ws.Range("H70").Value = 5708479.5
$ws.Range("I70").Value = 2504317.2
$ws.Range("J70").Value = 12990667
$ws.Range("K70").Value = 2504317.2
$ws.Range("L70").Value = 12990667
$ws.Range("M70").Value = -2504047.2
$ws.Range("N70").Value = -12991207

$ws.Range("H73").Value = 5708479.5
$ws.Range("I73").Value = 2504317.2
$ws.Range("J73").Value = 12990667
$ws.Range("K73").Value = 2504317.2
$ws.Range("L73").Value = 12990667
$ws.Range("M73").Value = -2503381.2
$ws.Range("N73").Value = -12992539

$ws.Range("H80").Value = 8896.429
$ws.Range("I80").Value = 4090
$ws.Range("J80").Value = 20912.5
$ws.Range("K80").Value = 4090
$ws.Range("L80").Value = 20912.5
$ws.Range("M80").Value = -3092
$ws.Range("N80").Value = -22908.5

$ws.Range("H83").Value = 8896.429
$ws.Range("I83").Value = 4090
$ws.Range("J83").Value = 20912.5
$ws.Range("K83").Value = 20450
$ws.Range("L83").Value = 104562.5
$ws.Range("M83").Value = -15458
$ws.Range("N83").Value = -114546.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 41667660
$ws.Range("I46").Value = 892.5714
$ws.Range("K46").Value = 892.5714
$ws.Range("M46").Value = -704.5714

$ws.Range("H93").Value = 12458.083
$ws.Range("I93").Value = 3064
$ws.Range("J93").Value = 25609.8
$ws.Range("K93").Value = 3064
$ws.Range("L93").Value = 25609.8
$ws.Range("M93").Value = -1816
$ws.Range("N93").Value = -28105.8

$ws.Range("H132").Value = 4613078.5
$ws.Range("I132").Value = 5295897
$ws.Range("J132").Value = 4051
$ws.Range("K132").Value = 15887691
$ws.Range("L132").Value = 12153
$ws.Range("M132").Value = -15885161
$ws.Range("N132").Value = -17213
